$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.131.65"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.909.33"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "348.41"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.16"
$ws.Range("E6").Value = "  -7.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.552"
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.42"
$ws.Range("E10").Value = "  -5.54%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.77"
$ws.Range("E13").Value = "  -6.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.365.08"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.902.42"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.952"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.138.03"
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.40"
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.29"
$ws.Range("E21").Value = "  -6.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.46"
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.53"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.68"
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.172"
$ws.Range("E26").Value = "  -4.50%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.19"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +4.63%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.16"
$ws.Range("E31").Value = "  -4.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.11"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.20"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.28"
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.18"
$ws.Range("E35").Value = "  -5.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -7.04%  "
$ws.Range("E38").Value = "  -8.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.51"
$ws.Range("E39").Value = "  -6.22%  "
$ws.Range("E40").Value = "  -6.32%  "
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.27"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.73"
$ws.Range("E44").Value = "  +6.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.090.08"
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  -6.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.26"
$ws.Range("E48").Value = "  -9.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.238"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("E50").Value = "  -4.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.882"
$ws.Range("E51").Value = "  -7.72%  "